$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 235, shifting existing rows 235:348 down to 236:349
$ws.Rows("235:235").Insert()

# Populate the newly inserted row 235 with the new record
$ws.Range("A235").Value = 10
$ws.Range("B235").Value = "Vega Modelo de Temuco"
$ws.Range("C235").Value = "La Araucanía"
$ws.Range("D235").Value = 44825
$ws.Range("E235").Value = 9
$ws.Range("F235").Value = 100112001
$ws.Range("G235").Value = "Berenjena"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 20
$ws.Range("K235").Value = 15000
$ws.Range("L235").Value = 15000
$ws.Range("M235").Value = 15000
$ws.Range("N235").Value = "$/caja 40 unidades"
$ws.Range("O235").Value = "Región de Arica y Parinacota"
$ws.Range("P235").Value = 375
$ws.Range("Q235").Value = 40
$ws.Range("R235").Value = "Hortaliza"
